{"js": "// The document header shows a date, followed by a 5-column practice table\n// of two-digit-by-two-digit multiplication equations (e.g. \"38x48=1824\").\n// This edit refreshes the date and swaps in a new set of equations, in the\n// same order the old ones appear in the document. Every old/new string below\n// is unique in the document, so we can safely find-and-replace each run's\n// text directly (this keeps each run's original formatting -- font/size --\n// untouched, since only the text content of the run is modified).\nconst replacements = [\n  [\"2025-06-11 Wednesday\", \"2025-06-12 Thursday\"],\n  [\"38\u00d748=1824\", \"16\u00d778=1248\"],\n  [\"68\u00d785=5780\", \"30\u00d759=1770\"],\n  [\"46\u00d755=2530\", \"74\u00d727=1998\"],\n  [\"16\u00d784=1344\", \"43\u00d780=3440\"],\n  [\"64\u00d754=3456\", \"72\u00d738=2736\"],\n  [\"66\u00d755=3630\", \"87\u00d747=4089\"],\n  [\"37\u00d761=2257\", \"17\u00d784=1428\"],\n  [\"86\u00d781=6966\", \"80\u00d784=6720\"],\n  [\"49\u00d737=1813\", \"55\u00d728=1540\"],\n  [\"79\u00d753=4187\", \"86\u00d787=7482\"],\n  [\"13\u00d748=624\", \"99\u00d762=6138\"],\n  [\"19\u00d761=1159\", \"38\u00d793=3534\"],\n  [\"22\u00d783=1826\", \"20\u00d768=1360\"],\n  [\"49\u00d766=3234\", \"54\u00d779=4266\"],\n  [\"26\u00d731=806\", \"29\u00d716=464\"],\n  [\"60\u00d747=2820\", \"85\u00d762=5270\"],\n  [\"94\u00d759=5546\", \"25\u00d736=900\"],\n  [\"18\u00d750=900\", \"22\u00d723=506\"],\n  [\"36\u00d796=3456\", \"46\u00d776=3496\"],\n  [\"35\u00d795=3325\", \"56\u00d768=3808\"],\n  [\"14\u00d784=1176\", \"31\u00d736=1116\"],\n  [\"51\u00d723=1173\", \"69\u00d780=5520\"],\n  [\"20\u00d735=700\", \"26\u00d740=1040\"],\n  [\"23\u00d731=713\", \"77\u00d757=4389\"],\n  [\"52\u00d791=4732\", \"54\u00d754=2916\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find text to replace: \" + oldText);\n  }\n\n  // Replace the matched range's text in place (Word.InsertLocation.replace)\n  // so the surrounding run formatting (rFonts/sz) is preserved.\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Apply the edit: update the header date and the 25 multiplication\n# equations in the practice-sheet table, in document order.\n# Each (old, new) pair is unique within the document, so Find/Replace\n# with MatchCase locates the exact run and substitutes its text in\n# place, leaving the run's formatting (font/size) untouched.\n\n$d = $word.ActiveDocument\n\n$wdFindContinue = 1\n$wdReplaceOne   = 1\n\n$replacements = @(\n    @('2025-06-11 Wednesday', '2025-06-12 Thursday'),\n    @('38\u00d748=1824', '16\u00d778=1248'),\n    @('68\u00d785=5780', '30\u00d759=1770'),\n    @('46\u00d755=2530', '74\u00d727=1998'),\n    @('16\u00d784=1344', '43\u00d780=3440'),\n    @('64\u00d754=3456', '72\u00d738=2736'),\n    @('66\u00d755=3630', '87\u00d747=4089'),\n    @('37\u00d761=2257', '17\u00d784=1428'),\n    @('86\u00d781=6966', '80\u00d784=6720'),\n    @('49\u00d737=1813', '55\u00d728=1540'),\n    @('79\u00d753=4187', '86\u00d787=7482'),\n    @('13\u00d748=624', '99\u00d762=6138'),\n    @('19\u00d761=1159', '38\u00d793=3534'),\n    @('22\u00d783=1826', '20\u00d768=1360'),\n    @('49\u00d766=3234', '54\u00d779=4266'),\n    @('26\u00d731=806', '29\u00d716=464'),\n    @('60\u00d747=2820', '85\u00d762=5270'),\n    @('94\u00d759=5546', '25\u00d736=900'),\n    @('18\u00d750=900', '22\u00d723=506'),\n    @('36\u00d796=3456', '46\u00d776=3496'),\n    @('35\u00d795=3325', '56\u00d768=3808'),\n    @('14\u00d784=1176', '31\u00d736=1116'),\n    @('51\u00d723=1173', '69\u00d780=5520'),\n    @('20\u00d735=700', '26\u00d740=1040'),\n    @('23\u00d731=713', '77\u00d757=4389'),\n    @('52\u00d791=4732', '54\u00d754=2916'),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n\n    $found = $find.Execute(\n        $oldText,          # FindText\n        $true,             # MatchCase\n        $false,            # MatchWholeWord\n        $false,            # MatchWildcards\n        $false,            # MatchSoundsLike\n        $false,            # MatchAllWordForms\n        $true,             # Forward\n        $wdFindContinue,   # Wrap\n        $false,            # Format\n        $newText,          # ReplaceWith\n        $wdReplaceOne      # Replace\n    )\n\n    if (-not $found) {\n        throw \"Could not find text to replace: $oldText\"\n    }\n}\n\n"}
